# Updates "cryptos" price-table worksheet with refreshed price/volume
# data (GitHub Actions scheduled refresh), including two pairs of rows
# whose coin identity + data swapped position (37/38 and 41/42 and 44/45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.437.81", "0.00001010")
# that must stay plain text -- force the text number format before
# assigning, then drop back to the default style so no stray cell
# style is left behind (matches original formatting exactly).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.205.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.440.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9157"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "275.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3074"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.89"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.024"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06492"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9988"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.345"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.050"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001010"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.436.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9354"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05624"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.391"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.238"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.204.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.136"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.586.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.832"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8107"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.833"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07635"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.467"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05842"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.656"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.128"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.03%  "
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1846"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9221"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.180"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -13.44%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5206"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.490"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5095"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.734"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06344"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9871"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
